# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font + border + centered alignment)
# from the last existing header cell (AC1) onto the three new header
# cells so they match the look of the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate every data row (2-47) with the team's season record.
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 89
    $ws.Cells.Item($row, 31).Value = 73
    $ws.Cells.Item($row, 32).Value = 0
}
